$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'259.85"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'5.98%"
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'28.01"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'-4.04%"
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'5.216"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'-0.65%"
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'0.05928"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'3.96%"
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'6.726"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'1.67%"
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.8743"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'2.72%"
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'1.008"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'18.11%"
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.1427"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'4.12%"
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.07247"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'2.52%"
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.03201"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'0.31%"
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.09250"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'0.17%"
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.001544"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'1.38%"
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'0.0006072"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'1.83%"
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.005879"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'-1.88%"
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'3.491"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'-0.13%"
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'1.27%"
$ws.Range('E17').Style = 'Normal'
$ws.Range('D19').Value = "'0.3124"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'-1.10%"
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.03653"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'12.66%"
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'0.1289"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'0.93%"
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'3.515"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'0.63%"
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.04180"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'2.36%"
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'1.26%"
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'-0.44%"
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'0.004570"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'10.38%"
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'0.0001198"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'-0.15%"
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'33.51%"
$ws.Range('E28').Style = 'Normal'
$ws.Range('D40').Value = "'0.03848"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'2.54%"
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.005433"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'46.36%"
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'4.09%"
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.002374"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'-1.20%"
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.01090"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'16.39%"
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.00005416"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'2.29%"
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.00000000749"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'-0.19%"
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.08539"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'13.79%"
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'0.002139"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'-12.40%"
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'0.00002097"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'-0.19%"
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.0001997"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'-0.19%"
$ws.Range('E50').Style = 'Normal'
